$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = "16:16"
$ws.Range("G2").Value = 200.9472680245446
$ws.Range("H2").Value = 99.93525306664431

# Row 3
$ws.Range("E3").Value = "16:16"
$ws.Range("F3").Value = [double]"4.035040301886283e-09"
$ws.Range("G3").Value = 397.8932912197835
$ws.Range("H3").Value = 85.10312446386273

# Row 4
$ws.Range("E4").Value = "16:16"
$ws.Range("F4").Value = 0.7327678547907854
$ws.Range("G4").Value = 998.3112709787323
$ws.Range("H4").Value = 69.99392331902639

# Row 5
$ws.Range("E5").Value = "16:16"
$ws.Range("F5").Value = 0.5913465965666089
$ws.Range("G5").Value = 400.5308692093085
$ws.Range("H5").Value = 54.90910915961167

# Row 6
$ws.Range("E6").Value = "16:16"
$ws.Range("F6").Value = 0.6392115640936014
$ws.Range("G6").Value = 998.8757817610992
$ws.Range("H6").Value = 130.0015809752132

# Row 7
$ws.Range("E7").Value = "16:16"
$ws.Range("F7").Value = 0.8295477617416493
$ws.Range("G7").Value = 1999.605449741923
$ws.Range("H7").Value = 114.9985694223957
